$wb = $excel.ActiveWorkbook

# New row (59) data to append to each of the 4 worksheets.
# Columns: A time, B total-len(hex), C ID(hex), D actual-len(hex), E checksum(hex),
#          F total-len_DEC (number), G ID_DEC (big-int stored as text), H actual-len_DEC (number), I checksum_DEC (number)
$rows = @{
    "ROW35-FE-LIFTER"  = @("2025-03-06 18:42:06", "0x01,0x90 ", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,", "0x01,0x90,", "0x d", 400, "568631262647113770877196", 400, 13)
    "ROW35-MID-LIFTER" = @("2025-03-06 18:29:35", "0x01,0x90 ", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,", "0x01,0x90,", "0x e", 400, "568631262647113770942732", 400, 14)
    "ROW02-FE-LIFTER"  = @("2025-03-06 18:51:45", "0x01,0x90 ", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,", "0x01,0x90,", "0xff",  400, "568631262647113769959692", 400, 255)
    "ROW02-MID-LIFTER" = @("2025-03-06 18:41:15", "0x01,0x90 ", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x01,0x90,", "0x 3",  400, "568631262647113769959692", 400, 3)
}

foreach ($sheetName in $rows.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $vals = $rows[$sheetName]
    $r = 59

    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
    $ws.Cells.Item($r, 6).Value = $vals[5]

    # Column G holds a 24-digit integer that exceeds double precision, so it
    # must stay text (matches the rest of the column). A leading apostrophe
    # forces text entry; reset the style afterwards so no stray number
    # format/style index gets attached to the cell.
    $ws.Cells.Item($r, 7).Value = "'" + $vals[6]
    $ws.Cells.Item($r, 7).Style = "Normal"

    $ws.Cells.Item($r, 8).Value = $vals[7]
    $ws.Cells.Item($r, 9).Value = $vals[8]
}
